# Generate Report for Handoff
#
# The localization report workbook tracks, per source file, the current
# hand-off/translation status for each target locale. The file
# "7e6bb69b-92a2-46b6-a799-327c287ce556.md" has just been packaged and is
# now "Ready for handoff" (it was previously "In Translation"). This script
# updates the Overview sheet as well as the per-locale (zh-cn / de-de)
# detail sheets to reflect the new handoff, refreshing the status, the
# translation priority, and the handoff timestamps.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet: row 3 corresponds to 7e6bb69b-92a2-46b6-a799-327c287ce556
# Columns: A=File Name, B=Path And Name, C=Extension, D=Publish URL,
#          E=zh-cn, F=de-de, G=Latest HO Xliff Generate Date
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-08-16 12:13:25"
$overview.Columns.Item(5).ColumnWidth = 16.3
$overview.Columns.Item(6).ColumnWidth = 16.3

# ---------------------------------------------------------------------
# zh-cn sheet: row 3 corresponds to 7e6bb69b-92a2-46b6-a799-327c287ce556
# Columns: A=Source File Name, B=File Extension, C=Status, D=Source Path,
#          E=Priority, F=Content Duplicate, G=Latest Handoff File,
#          H=Latest Handoff Datetime, ...
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("E3").Value = "mt"
$zhcn.Range("G3").Value = "7e6bb69b-92a2-46b6-a799-327c287ce556.92a671c022230d45ac5a12e1d798625ae497d00b.zh-cn.xlf"
$zhcn.Range("H3").Value = "2016-08-16 12:13:21"
$zhcn.Columns.Item(3).ColumnWidth = 16.3

# ---------------------------------------------------------------------
# de-de sheet: row 3 corresponds to 7e6bb69b-92a2-46b6-a799-327c287ce556
# Columns: A=Source File Name, B=File Extension, C=Status, D=Source Path,
#          E=Priority, F=Content Duplicate, G=Latest Handoff File,
#          H=Latest Handoff Datetime, ...
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("E3").Value = "mt"
$dede.Range("G3").Value = "7e6bb69b-92a2-46b6-a799-327c287ce556.92a671c022230d45ac5a12e1d798625ae497d00b.de-de.xlf"
$dede.Range("H3").Value = "2016-08-16 12:13:25"
$dede.Columns.Item(3).ColumnWidth = 16.3
